$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.048.98"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "2.172.46"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.97"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.94%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0927"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.100"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "2.498.23"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.809"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.05%  "
$ws.Range("D17").Value = "2.161.50"
$ws.Range("E17").Value = "  -4.46%  "
$ws.Range("D18").Value = "40.972.01"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("E19").Value = "  -7.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.65%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.88%  "
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.41%  "
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.105"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0286"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "60.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.190"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0972"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("E49").Value = "  -6.63%  "
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("D51").Value = "2.373.21"
$ws.Range("E51").Value = "  -2.13%  "
